# Update countries & provincias Spain
# - Refresh COVID-19 case figures for a batch of countries
# - Maldivas overtakes Paraguay in total cases -> rows swap order
# - Islas Malvinas / Groenlandia swap order (tied totals)
# - Refresh "last updated" timestamp

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Timestamp header (row 1) ---
$ws.Cells.Item(1, 1).Value = "Datos actualizados a 8 de Julio de 2020 a las 20:27"

# --- Simple per-country numeric refreshes (country/row mapping unchanged) ---

# Row 4: Estados Unidos
$ws.Cells.Item(4, 2).Value = 3122040
$ws.Cells.Item(4, 3).Value = 24956
$ws.Cells.Item(4, 4).Value = 1372830
$ws.Cells.Item(4, 5).Value = 1614863
$ws.Cells.Item(4, 7).Value = 375
$ws.Cells.Item(4, 8).Value = 134347

# Row 6: India
$ws.Cells.Item(6, 2).Value = 769041
$ws.Cells.Item(6, 3).Value = 25560
$ws.Cells.Item(6, 4).Value = 476550
$ws.Cells.Item(6, 5).Value = 271347
$ws.Cells.Item(6, 7).Value = 491
$ws.Cells.Item(6, 8).Value = 21144

# Row 19: Alemania
$ws.Cells.Item(19, 2).Value = 198703
$ws.Cells.Item(19, 3).Value = 348
$ws.Cells.Item(19, 5).Value = 6894
$ws.Cells.Item(19, 7).Value = 6
$ws.Cells.Item(19, 8).Value = 9109

# Row 21: Francia
$ws.Cells.Item(21, 2).Value = 169473
$ws.Cells.Item(21, 3).Value = 663
$ws.Cells.Item(21, 5).Value = 61853
$ws.Cells.Item(21, 7).Value = 32
$ws.Cells.Item(21, 8).Value = 29965

# Row 36: Kazajistan
$ws.Cells.Item(36, 4).Value = 34149
$ws.Cells.Item(36, 5).Value = 16646

# Row 66: Marruecos
$ws.Cells.Item(66, 2).Value = 14771
$ws.Cells.Item(66, 3).Value = 164
$ws.Cells.Item(66, 4).Value = 11316
$ws.Cells.Item(66, 5).Value = 3213
$ws.Cells.Item(66, 7).Value = 2
$ws.Cells.Item(66, 8).Value = 242

# Row 72: Sudan
$ws.Cells.Item(72, 2).Value = 10084
$ws.Cells.Item(72, 3).Value = 87
$ws.Cells.Item(72, 4).Value = 5074
$ws.Cells.Item(72, 5).Value = 4374
$ws.Cells.Item(72, 7).Value = 14
$ws.Cells.Item(72, 8).Value = 636

# Row 95: Republica de Yibuti
$ws.Cells.Item(95, 2).Value = 4889
$ws.Cells.Item(95, 3).Value = 11
$ws.Cells.Item(95, 4).Value = 4644
$ws.Cells.Item(95, 5).Value = 190

# Row 98: Republica de Africa Central
$ws.Cells.Item(98, 2).Value = 4109
$ws.Cells.Item(98, 3).Value = 38
$ws.Cells.Item(98, 4).Value = 1050
$ws.Cells.Item(98, 5).Value = 3007

# Row 105: Somalia
$ws.Cells.Item(105, 2).Value = 3028
$ws.Cells.Item(105, 3).Value = 13
$ws.Cells.Item(105, 4).Value = 1147
$ws.Cells.Item(105, 5).Value = 1789

# Row 128: Yemen
$ws.Cells.Item(128, 2).Value = 1318
$ws.Cells.Item(128, 3).Value = 21
$ws.Cells.Item(128, 4).Value = 595
$ws.Cells.Item(128, 5).Value = 372
$ws.Cells.Item(128, 7).Value = 3
$ws.Cells.Item(128, 8).Value = 351

# Row 133: Jordania
$ws.Cells.Item(133, 4).Value = 977
$ws.Cells.Item(133, 5).Value = 182

# Row 164: Birmania
$ws.Cells.Item(164, 2).Value = 317
$ws.Cells.Item(164, 3).Value = 1
$ws.Cells.Item(164, 4).Value = 250
$ws.Cells.Item(164, 5).Value = 61

# --- Maldivas overtakes Paraguay: rows 108/109 swap country + data ---

# Row 108 becomes Maldivas (updated figures)
$ws.Cells.Item(108, 1).Value = "Maldivas"
$ws.Cells.Item(108, 2).Value = 2517
$ws.Cells.Item(108, 3).Value = 16
$ws.Cells.Item(108, 4).Value = 2180
$ws.Cells.Item(108, 5).Value = 324
$ws.Cells.Item(108, 6).Value = 0
$ws.Cells.Item(108, 7).Value = 1
$ws.Cells.Item(108, 8).Value = 13

# Row 109 becomes Paraguay (figures unchanged from before)
$ws.Cells.Item(109, 1).Value = "Paraguay"
$ws.Cells.Item(109, 2).Value = 2502
$ws.Cells.Item(109, 3).Value = 0
$ws.Cells.Item(109, 4).Value = 1193
$ws.Cells.Item(109, 5).Value = 1289
$ws.Cells.Item(109, 6).Value = 0
$ws.Cells.Item(109, 7).Value = 0
$ws.Cells.Item(109, 8).Value = 20

# --- Islas Malvinas / Groenlandia swap order (tied totals, no numeric change) ---

$ws.Cells.Item(209, 1).Value = "Islas Malvinas"
$ws.Cells.Item(210, 1).Value = "Groenlandia"
